$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3911.2222
$ws.Range("I70").Value = 10002
$ws.Range("J70").Value = 3149.875
$ws.Range("K70").Value = 30006
$ws.Range("L70").Value = 9449.625
$ws.Range("M70").Value = -29736
$ws.Range("N70").Value = -9989.625

$ws.Range("H73").Value = 3911.2222
$ws.Range("I73").Value = 10002
$ws.Range("J73").Value = 3149.875
$ws.Range("K73").Value = 30006
$ws.Range("L73").Value = 9449.625
$ws.Range("M73").Value = -29070
$ws.Range("N73").Value = -11321.625

$ws.Range("H107").Value = 453.05264
$ws.Range("I107").Value = 476.53333
$ws.Range("J107").Value = 365
$ws.Range("K107").Value = 476.53333
$ws.Range("L107").Value = 365
$ws.Range("M107").Value = 1443.46667
$ws.Range("N107").Value = -4205

$ws.Range("H111").Value = 250422.5
$ws.Range("I111").Value = 563
$ws.Range("J111").Value = 1000001
$ws.Range("K111").Value = 1689
$ws.Range("L111").Value = 3000003
$ws.Range("M111").Value = 1378
$ws.Range("N111").Value = -3006137

$ws.Range("H113").Value = 2778966.2
$ws.Range("I113").Value = 5556610
$ws.Range("J113").Value = 1322.6666
$ws.Range("K113").Value = 5556610
$ws.Range("L113").Value = 1322.6666
$ws.Range("M113").Value = -5553356
$ws.Range("N113").Value = -7830.6666

$ws.Range("H132").Value = 1981.8276
$ws.Range("I132").Value = 1556.1923
$ws.Range("K132").Value = 4668.5769
$ws.Range("M132").Value = -2138.5769


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 273393.78
$ws.Range("I61").Value = 2446.1667
$ws.Range("J61").Value = 530081
$ws.Range("K61").Value = 2446.1667
$ws.Range("L61").Value = 530081
$ws.Range("M61").Value = -2234.1667
$ws.Range("N61").Value = -530505

$ws.Range("H74").Value = 1676.6666
$ws.Range("I74").Value = 1337.8096
$ws.Range("K74").Value = 1337.8096
$ws.Range("M74").Value = -463.8096

$ws.Range("H77").Value = 1676.6666
$ws.Range("I77").Value = 1337.8096
$ws.Range("K77").Value = 6689.048000000001
$ws.Range("M77").Value = -2321.048000000001

$ws.Range("H110").Value = 1262.8
$ws.Range("I110").Value = 1009.05
$ws.Range("J110").Value = 2277.8
$ws.Range("K110").Value = 1009.05
$ws.Range("L110").Value = 2277.8
$ws.Range("M110").Value = 1035.95
$ws.Range("N110").Value = -6367.8

$ws.Range("H122").Value = 1639.091
$ws.Range("I122").Value = 1673.0322
$ws.Range("J122").Value = 1113
$ws.Range("K122").Value = 5019.096600000001
$ws.Range("L122").Value = 3339
$ws.Range("M122").Value = -2569.096600000001
$ws.Range("N122").Value = -8239

$ws.Range("H136").Value = 273393.78
$ws.Range("I136").Value = 2446.1667
$ws.Range("J136").Value = 530081
$ws.Range("K136").Value = 7338.500100000001
$ws.Range("L136").Value = 1590243
$ws.Range("M136").Value = -4788.500100000001
$ws.Range("N136").Value = -1595343


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 941.29034
$ws.Range("I107").Value = 786.2917
$ws.Range("J107").Value = 1472.7142
$ws.Range("K107").Value = 786.2917
$ws.Range("L107").Value = 1472.7142
$ws.Range("M107").Value = 1133.7083
$ws.Range("N107").Value = -5312.7142


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3373.56
$ws.Range("I31").Value = 2467.8462
$ws.Range("J31").Value = 3691.7837
$ws.Range("K31").Value = 2467.8462
$ws.Range("L31").Value = 3691.7837
$ws.Range("M31").Value = -2172.8462
$ws.Range("N31").Value = -4281.7837

$ws.Range("H34").Value = 3373.56
$ws.Range("I34").Value = 2467.8462
$ws.Range("J34").Value = 3691.7837
$ws.Range("K34").Value = 2467.8462
$ws.Range("L34").Value = 3691.7837
$ws.Range("M34").Value = -2265.8462
$ws.Range("N34").Value = -4095.7837

$ws.Range("H107").Value = 532.8148
$ws.Range("I107").Value = 211.83333
$ws.Range("J107").Value = 1174.7778
$ws.Range("K107").Value = 211.83333
$ws.Range("L107").Value = 1174.7778
$ws.Range("M107").Value = 1708.16667
$ws.Range("N107").Value = -5014.7778

$ws.Range("H118").Value = 39464.652
$ws.Range("J118").Value = 39464.652
$ws.Range("L118").Value = 39464.652
$ws.Range("N118").Value = -42778.652

$ws.Range("H132").Value = 2425.8076
$ws.Range("I132").Value = 2050.2942
$ws.Range("J132").Value = 3135.111
$ws.Range("K132").Value = 6150.882599999999
$ws.Range("L132").Value = 9405.332999999999
$ws.Range("M132").Value = -3620.882599999999
$ws.Range("N132").Value = -14465.333

$ws.Range("H134").Value = 4075.4
$ws.Range("I134").Value = 4105.684
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 12317.052
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -9782.052
$ws.Range("N134").Value = -15570


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1325173.6
$ws.Range("I122").Value = 2850604.8
$ws.Range("K122").Value = 8551814.399999999
$ws.Range("M122").Value = -8549364.399999999

$ws.Range("H123").Value = 25099.55
$ws.Range("J123").Value = 25946.895
$ws.Range("L123").Value = 25946.895
$ws.Range("N123").Value = -30846.895


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 360700.6
$ws.Range("I122").Value = 3852.3157
$ws.Range("J122").Value = 1114047
$ws.Range("K122").Value = 11556.9471
$ws.Range("L122").Value = 3342141
$ws.Range("M122").Value = -9106.947100000001
$ws.Range("N122").Value = -3347041

$ws.Range("H136").Value = 4845.564
$ws.Range("I136").Value = 2283.3794
$ws.Range("J136").Value = 12275.9
$ws.Range("K136").Value = 6850.138199999999
$ws.Range("L136").Value = 36827.7
$ws.Range("M136").Value = -4300.138199999999
$ws.Range("N136").Value = -41927.7


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4733.3335
$ws.Range("I62").Value = 4225
$ws.Range("J62").Value = 5750
$ws.Range("K62").Value = 4225
$ws.Range("L62").Value = 5750
$ws.Range("M62").Value = -3601
$ws.Range("N62").Value = -6998

$ws.Range("H65").Value = 4733.3335
$ws.Range("I65").Value = 4225
$ws.Range("J65").Value = 5750
$ws.Range("K65").Value = 21125
$ws.Range("L65").Value = 28750
$ws.Range("M65").Value = -18005
$ws.Range("N65").Value = -34990

$ws.Range("H113").Value = 661.1786
$ws.Range("I113").Value = 484.85
$ws.Range("J113").Value = 1102
$ws.Range("K113").Value = 1454.55
$ws.Range("L113").Value = 3306
$ws.Range("M113").Value = 715.4499999999998
$ws.Range("N113").Value = -7646

$ws.Range("H132").Value = 1654.3948
$ws.Range("I132").Value = 1246.3182
$ws.Range("J132").Value = 2215.5
$ws.Range("K132").Value = 3738.9546
$ws.Range("L132").Value = 6646.5
$ws.Range("M132").Value = -1208.9546
$ws.Range("N132").Value = -11706.5

